$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated T3 DRC data (columns C "T4" and D "T5") for rows 2-6
$ws.Range("C2").Value = 37.5513798303019
$ws.Range("D2").Value = 37.41191392868211

$ws.Range("C3").Value = 38.19481299091082
$ws.Range("D3").Value = 38.8683304223175

$ws.Range("C4").Value = 37.93527167106304
$ws.Range("D4").Value = 38.12493431171944

$ws.Range("C5").Value = 37.30760049005133
$ws.Range("D5").Value = 37.74486609977946

$ws.Range("C6").Value = 37.45934422616917
$ws.Range("D6").Value = 37.3459741174844
